$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Split the old "extension" sheet into two sheets:
#    - Rename the existing "extension" sheet (sheetId 29) to "AnyValue" and
#      clear its contents (AnyValue has no header row at all).
#    - Add a brand new sheet right after it, named "extension", which gets
#      the extension_tag / extension_value / extensions header row that used
#      to live on the original sheet.
# ---------------------------------------------------------------------------
$oldExtension = $wb.Worksheets.Item("extension")
$oldExtension.Cells.Clear()
$oldExtension.Name = "AnyValue"

$newExtension = $wb.Worksheets.Add($null, $oldExtension)
$newExtension.Name = "extension"
$newExtension.Range("A1").Value = "extension_tag"
$newExtension.Range("B1").Value = "extension_value"
$newExtension.Range("C1").Value = "extensions"

# ---------------------------------------------------------------------------
# 2. Insert a new "id_prefixes_are_closed" column immediately after the
#    "id_prefixes" column on every sheet that has it (the *_definition
#    sheets).
# ---------------------------------------------------------------------------
$sheetsWithIdPrefixes = @(
    "schema_definition",
    "type_definition",
    "subset_definition",
    "enum_definition",
    "slot_definition",
    "class_definition"
)

foreach ($sheetName in $sheetsWithIdPrefixes) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the column that currently holds "id_prefixes" in row 1.
    $found = $ws.Rows.Item(1).Find("id_prefixes")
    $col = $found.Column

    $insertCol = $ws.Cells.Item(1, $col + 1).EntireColumn
    $insertCol.Insert(-4161)
    $ws.Cells.Item(1, $col + 1).Value = "id_prefixes_are_closed"
}

# ---------------------------------------------------------------------------
# 3. Extend the CODE/CURIE/URI/FHIR_CODING dropdown validations with LABEL.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("enum_expression")
$ws6.Range("D2:D1048576").Validation.Formula1 = '"CODE,CURIE,URI,FHIR_CODING,LABEL"'

$ws7 = $wb.Worksheets.Item("anonymous_enum_expression")
$ws7.Range("D2:D1048576").Validation.Formula1 = '"CODE,CURIE,URI,FHIR_CODING,LABEL"'

$ws8 = $wb.Worksheets.Item("enum_definition")
$ws8.Range("E2:E1048576").Validation.Formula1 = '"CODE,CURIE,URI,FHIR_CODING,LABEL"'
